# Merge - Opp Test Data, ENg Detail, Add Counterparty - 10 Oct 2025
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# AddOpportunity sheet: update the Sector value in the sample data row and
# normalize the stray formatting that had been applied to D2.
# ---------------------------------------------------------------------------
$wsAdd = $wb.Worksheets.Item("AddOpportunity")
$wsAdd.Range("C2").Value = "Debt Financing"
$wsAdd.Range("D2").Style = "Normal"

# ---------------------------------------------------------------------------
# Update the remembered cell selection on the Opportunity sheet.
# ---------------------------------------------------------------------------
$wsOpp = $wb.Worksheets.Item("Opportunity")
$wsOpp.Activate() | Out-Null
$wsOpp.Range("H19").Select() | Out-Null

# ---------------------------------------------------------------------------
# Update the remembered cell selection on the AddOpportunity sheet, and make
# sure it ends up as the active/visible sheet again (as it was originally).
# ---------------------------------------------------------------------------
$wsAdd.Activate() | Out-Null
$wsAdd.Range("E16").Select() | Out-Null
